$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Parneet kaur"

# Row 7 - __init__ / Attribute set to input values
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'color="red", length=5, width=6'
$ws.Range("G7").Value = "Rectangle object created successfully; attributes set to input values"

# Row 8 - __init__ / Exception raised when color is blank
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'ValueError("Color cannot be blank.")'
$ws.Range("G8").Value = 'ValueError("Color cannot be blank.")'

# Row 9 - __init__ / Exception raised when length is not an integer
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'color="red", length=5.7, width=6'
$ws.Range("G9").Value = 'ValueError("Length must be numeric.")'

# Row 10 - __init__ / Exception raised when width is not an integer
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'color="red", length=5, width="6"'
$ws.Range("G10").Value = 'ValueError("Width must be numeric.")'

# Row 11 - __str__ / Returns string formatted appropriately
$ws.Range("E11").Value = 'Valid rectangle exists (e.g., Rectangle("red", 5, 6))'
$ws.Range("F11").Value = "Call str(rectangle)"
$ws.Range("G11").Value = "String contains both`nThe shape color is red.`n5, 6, 5 and 6"

# Row 12 - calculate_area / Returns correct calculated value
$ws.Range("E12").Value = 'Valid rectangle exists (e.g., Rectangle("blue", 3, 4))'
$ws.Range("F12").Value = "Call calculate_area()"
$ws.Range("G12").Value = "Returns 12.0"

# Row 13 - calculate_perimeter / Returns correct calculated value
$ws.Range("E13").Value = 'Valid rectangle exists (e.g., Rectangle("green", 2, 3))'
$ws.Range("F13").Value = "Call calculate_perimeter()"
$ws.Range("G13").Value = "Returns 10.0"

# Update selection / scroll position to match saved view state
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 7
[void]$ws.Range("G13").Select()
